$wb = $excel.ActiveWorkbook

$schedule = $wb.Worksheets.Item("Schedule")
$detailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates ---
$schedule.Range("E3").Value = 278.1123345
$schedule.Range("F3").Value = 18.39367291666667
$schedule.Range("E4").Value = -131.59852875
$schedule.Range("F4").Value = -2.901202132936508

# --- Detailed sheet updates ---
$detailed.Range("B35").Value = 8.64973
$detailed.Range("B36").Value = -7.981
$detailed.Range("B37").Value = -8.90727
$detailed.Range("C37").Value = "historical"
$detailed.Range("B38").Value = -2.98843
$detailed.Range("C38").Value = "historical"
$detailed.Range("B39").Value = -0.62019
$detailed.Range("B40").Value = 3.15941
$detailed.Range("B41").Value = 24.11812
$detailed.Range("B42").Value = 25.00316
$detailed.Range("B43").Value = 17.20654
$detailed.Range("B44").Value = 16.67999
$detailed.Range("B45").Value = 56.98
$detailed.Range("B46").Value = 36.0601
$detailed.Range("B47").Value = 53.84201
$detailed.Range("B48").Value = 43.41162
$detailed.Range("B49").Value = 36.06
$detailed.Range("B52").Value = 40.54
$detailed.Range("B54").Value = 35.87992
$detailed.Range("B59").Value = 65
$detailed.Range("B60").Value = 57.1
$detailed.Range("B61").Value = 57.1
$detailed.Range("B64").Value = 36.05971
$detailed.Range("B66").Value = -0.87734
$detailed.Range("B67").Value = -5.65164
$detailed.Range("B68").Value = -5.30295
$detailed.Range("B69").Value = -5.50985
$detailed.Range("B70").Value = -0.87893
$detailed.Range("B71").Value = -4.64736
$detailed.Range("B72").Value = -5.01
$detailed.Range("B73").Value = -1.092
$detailed.Range("B74").Value = -5.06248
$detailed.Range("B75").Value = -5.63691
$detailed.Range("B76").Value = -6.8
$detailed.Range("B77").Value = -11.16992
$detailed.Range("B78").Value = -12.35725
$detailed.Range("B81").Value = -8.222300000000001
$detailed.Range("B82").Value = 0
$detailed.Range("B83").Value = -7.74498
$detailed.Range("B84").Value = -7.78104
$detailed.Range("B85").Value = -5.66127
$detailed.Range("B86").Value = -6.18295
$detailed.Range("B87").Value = -6.02102
$detailed.Range("B90").Value = 43.04957
$detailed.Range("B91").Value = 56.98
$detailed.Range("B92").Value = 8.320819999999999
$detailed.Range("B97").Value = 65
